$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 3).Value = 'Alty looks at the ''key'' in her hand.' + "`n" + ''
$ws.Cells.Item(116, 3).Value = '[name="Frost"]  I’ll call this song... ''D''!' + "`n" + ''
$ws.Cells.Item(132, 3).Value = '[name="Archosauria B"]  It might be like what the priests do, it’s called ''music''?' + "`n" + ''
$ws.Cells.Item(134, 3).Value = '[name="Archosauria D"]  I know, they must be ''Cuicayotl''!' + "`n" + ''
$ws.Cells.Item(135, 3).Value = '[name="Dan"]  ''Cuicayotl''?  What’s that?' + "`n" + ''
$ws.Cells.Item(136, 3).Value = '[name="High Priest"]  Oh, that’s how you say ''one who sings'' in their language.' + "`n" + ''
$ws.Cells.Item(137, 3).Value = '[name="High Priest"]  But it’s been a long time since they had a ''Cuicayotl'' here, and you’ve won them over with your music.' + "`n" + ''
$ws.Cells.Item(138, 3).Value = '[name="High Priest"]  And I must say, your music is vastly different from that of any other ''Cuicayotl'' I’ve ever heard before!' + "`n" + ''
$ws.Cells.Item(140, 3).Value = '[name="Archosauria E"]  ''Cuicayotl'', give us another song!' + "`n" + ''
$ws.Cells.Item(156, 3).Value = '[name="Inam"]  That’s right. This is where all the tribes used to gather for the ''Mahuizzotia'', but it hasn’t been used since Gavial left.' + "`n" + ''
$ws.Cells.Item(157, 3).Value = '[name="Aya"]  ''Mahuizzotia''?' + "`n" + ''
$ws.Cells.Item(174, 3).Value = '[name="Inam"]  ''We left. Don’t miss us. I hope you enjoyed the music. Also, we left all our records here for you. Have fun.  -AUS''' + "`n" + ''
$ws.Cells.Item(196, 3).Value = 'Alty looks at the key in her hand, ''Dr. Kal''tsit, can that Rhodes Island of yours give the world the answer it’s looking for?''' + "`n" + ''
